$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(0.5692387035740197, 0.2355707868591003, 0.8033234571376444, 0.1059330575042835, 0.5682091107924527, 0.2557185134775685, 0.4537982025430372, 0.1529345173125635, 0.4618949180157207, 0.3074147176641421, 0.2823828768330753, 0.5056861808251918, 0.06015717143422494, 0.5272142884576241, 28.72735599237364, 44.57274171566024)

for ($row = 2; $row -le 26; $row++) {
    for ($col = 2; $col -le 17; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 2]
    }
}
